# Reorders the data rows (2..85) of the "Avverkningsanmälningar" sheet
# according to the permutation observed between the before/after OOXML,
# and bumps every "Förändrad" (column C) timestamp from 46063 to 46064.
#
# Strategy: snapshot every cell (value or formula) for rows 2..85 across
# columns A..Z, then write each target row's content back from the
# snapshot of its mapped source row. Finally force column C to 46064 for
# every data row (matches every row in the diff, including the ones whose
# row content didn't otherwise move).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 85
$lastCol = 26  # column Z

# target row -> source row (content copied from the source row's current
# position into the target row's position)
$rowMap = @{
    2=2; 3=6; 4=3; 5=5; 6=4; 7=7; 8=8; 9=9; 10=10; 11=11; 12=12; 13=13;
    14=14; 15=15; 16=16; 17=17; 18=18; 19=19; 20=20; 21=21; 22=22; 23=23;
    24=24; 25=25; 26=80; 27=74; 28=39; 29=73; 30=69; 31=43; 32=78; 33=84;
    34=66; 35=85; 36=77; 37=40; 38=41; 39=36; 40=44; 41=70; 42=79; 43=48;
    44=38; 45=71; 46=27; 47=61; 48=68; 49=28; 50=29; 51=51; 52=50; 53=58;
    54=55; 55=56; 56=54; 57=57; 58=53; 59=60; 60=59; 61=75; 62=82; 63=64;
    64=83; 65=45; 66=65; 67=33; 68=34; 69=46; 70=32; 71=47; 72=52; 73=30;
    74=37; 75=81; 76=35; 77=76; 78=49; 79=26; 80=63; 81=62; 82=31; 83=67;
    84=72; 85=42
}

# ---- 1. Snapshot every row currently on the sheet ----
$snapshot = @{}
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $rowData = @{}
    for ($col = 1; $col -le $lastCol; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        if ($cell.HasFormula) {
            $rowData[$col] = @{ IsFormula = $true; Val = $cell.Formula }
        } else {
            $rowData[$col] = @{ IsFormula = $false; Val = $cell.Value2 }
        }
    }
    $snapshot[$row] = $rowData
}

# ---- 2. Write each target row back from its mapped source row ----
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $srcRow = $rowMap[$row]
    $rowData = $snapshot[$srcRow]
    for ($col = 1; $col -le $lastCol; $col++) {
        $entry = $rowData[$col]
        $cell = $ws.Cells.Item($row, $col)
        if ($entry.IsFormula) {
            $cell.Formula = $entry.Val
        } else {
            $cell.Value = $entry.Val
        }
    }
}

# ---- 3. Every data row's "Förändrad" date moves from 46063 to 46064 ----
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 3).Value = 46064
}
